# Adding to work excel
# Appends a new tracking entry (row 91) to the "Sheet1" work log:
#   Item      = Refactoring of CardDecks to remove redundant abstract base
#               class, generify deck
#   Type      = Refactor
#   Entered   = 15-Dec-2018 (Excel serial 43449)
#   Assigned  = Andrew
#   Completed = 15-Dec-2018 (Excel serial 43449)
#   Outcome   = Completed, unit and integration tests all passing

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 91

$ws.Cells.Item($newRow, 1).Value = "Refactoring of CardDecks to remove redundant abstract base class, generify deck"
$ws.Cells.Item($newRow, 2).Value = "Refactor"

# Entered date - use the same short-date display ("d-mmm") as the rest of
# column C so the new row matches the existing formatting.
$ws.Cells.Item($newRow, 3).Value = 43449
$ws.Cells.Item($newRow, 3).NumberFormat = "d-mmm"

$ws.Cells.Item($newRow, 4).Value = "Andrew"

# Completed date - same treatment as column C / E elsewhere in the sheet.
$ws.Cells.Item($newRow, 5).Value = 43449
$ws.Cells.Item($newRow, 5).NumberFormat = "d-mmm"

$ws.Cells.Item($newRow, 6).Value = "Completed, unit and integration tests all passing"

# Leave the selection on the cell just past the newly-entered row, matching
# where the cursor lands after typing the last value in a new row.
$ws.Range("F92").Select()
